# "constrain and conflict check"
#
# The course/professor timetable on Sheet1 is reworked for an (external)
# scheduling conflict check: the old "Slots" column (C), which recorded a
# manually-assigned timetable slot per course, is dropped and replaced with
# a "Faculty" column (D) that mirrors the "Prof" column - the key the
# conflict checker constrains/joins on. A few course codes are also
# renumbered to their current catalogue IDs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Renumber course codes in column A (rows 2-7) ----------------------
$ws.Range("A2").Value = "IC218"
$ws.Range("A3").Value = "IC219"
$ws.Range("A4").Value = "IC222"
$ws.Range("A5").Value = "IC220"
$ws.Range("A6").Value = "IC223"
$ws.Range("A7").Value = "IC221"

# --- Drop the old "Slots" column (C) entirely ---------------------------
$ws.Range("C1:C7").ClearContents()

# --- Add the new "Faculty" column (D), mirroring column B ---------------
$ws.Range("D1").Value = "Faculty"
$ws.Range("D2").Value = $ws.Range("B2").Text
$ws.Range("D3").Value = $ws.Range("B3").Text
$ws.Range("D4").Value = $ws.Range("B4").Text
$ws.Range("D5").Value = $ws.Range("B5").Text
$ws.Range("D6").Value = $ws.Range("B6").Text
$ws.Range("D7").Value = $ws.Range("B7").Text

# --- Leave the sheet focused on the new Faculty column, as in the saved
#     selection state of the edited workbook.
$ws.Range("D1:E7").Select()
